$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 17:36"

# Country name swaps (data source re-sorted these pairs between the two updates)
$ws.Range("A102").Value = "Grecia"
$ws.Range("A103").Value = "Paraguay"
$ws.Range("A166").Value = "Mauricio"
$ws.Range("A167").Value = "Birmania"
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 4193667
$ws.Range("C4").Value = 23349
$ws.Range("D4").Value = 1981538
$ws.Range("E4").Value = 2064466
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 147663
$ws.Range("B6").Value = 1319302
$ws.Range("C6").Value = 31172
$ws.Range("D6").Value = 837459
$ws.Range("E6").Value = 450819
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 379
$ws.Range("H6").Value = 31024
$ws.Range("B13").Value = 297914
$ws.Range("C13").Value = 768
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 123
$ws.Range("H13").Value = 45677
$ws.Range("B17").Value = 245590
$ws.Range("C17").Value = 252
$ws.Range("D17").Value = 198192
$ws.Range("E17").Value = 12301
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 35097
$ws.Range("B21").Value = 205392
$ws.Range("C21").Value = 250
$ws.Range("D21").Value = 189400
$ws.Range("E21").Value = 6802
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 9190
$ws.Range("B23").Value = 148027
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 65447
$ws.Range("E23").Value = 79858
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 2722
$ws.Range("B24").Value = 112867
$ws.Range("C24").Value = 195
$ws.Range("D24").Value = 98519
$ws.Range("E24").Value = 5471
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 8877
$ws.Range("B44").Value = 52595
$ws.Range("C44").Value = 191
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 6139
$ws.Range("B63").Value = 22483
$ws.Range("C63").Value = 378
$ws.Range("D63").Value = 15407
$ws.Range("E63").Value = 6350
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 726
$ws.Range("B100").Value = 4570
$ws.Range("C100").Value = 104
$ws.Range("D100").Value = 2608
$ws.Range("E100").Value = 1834
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 5
$ws.Range("H100").Value = 128
$ws.Range("B102").Value = 4135
$ws.Range("C102").Value = 25
$ws.Range("D102").Value = 1374
$ws.Range("E102").Value = 2560
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 201
$ws.Range("B103").Value = 4113
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 2487
$ws.Range("E103").Value = 1590
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 36
$ws.Range("B116").Value = 2503
$ws.Range("C116").Value = 9
$ws.Range("D116").Value = 1901
$ws.Range("E116").Value = 479
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 123
$ws.Range("B147").Value = 1047
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 847
$ws.Range("E147").Value = 181
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 19
$ws.Range("B155").Value = 699
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 657
$ws.Range("E155").Value = 0
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 42
$ws.Range("B159").Value = 608
$ws.Range("C159").Value = 24
$ws.Range("D159").Value = 184
$ws.Range("E159").Value = 389
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 35
$ws.Range("B166").Value = 344
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 332
$ws.Range("E166").Value = 2
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 10
$ws.Range("B167").Value = 343
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 282
$ws.Range("E167").Value = 55
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 6
